$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 8-15 (C/D/E columns changed per diff) ---
# Row 8
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

# Row 9
$ws.Cells.Item(9, 3).Value = 16

# Row 10
$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 12
$ws.Cells.Item(10, 5).Value = $true

# Row 11
$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 9
$ws.Cells.Item(11, 5).Value = $true

# Row 12
$ws.Cells.Item(12, 3).Value = 10

# Row 13
$ws.Cells.Item(13, 4).Value = 8
$ws.Cells.Item(13, 5).Value = $true

# Row 14
$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 11
$ws.Cells.Item(14, 5).Value = $false

# Row 15
$ws.Cells.Item(15, 3).Value = 7
$ws.Cells.Item(15, 4).Value = 11
$ws.Cells.Item(15, 5).Value = $false

# --- Add new rows 16 and 17 ---
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "extr7"
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $false

$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "extr8"
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $true

# Match the style of the A column for the new rows (bold / bordered / centered)
foreach ($addr in @("A16", "A17")) {
    $r = $ws.Range($addr)
    $r.Font.Bold = $true
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4160
    $r.Borders.Item(7).LineStyle = 1
    $r.Borders.Item(10).LineStyle = 1
    $r.Borders.Item(8).LineStyle = 1
    $r.Borders.Item(9).LineStyle = 1
}
